$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Name values (A6, A7) - remove trailing F/R from URA3F/URA3R
$ws.Range("A6").Value = "YNL268W_sgtF_URA3"
$ws.Range("A7").Value = "YNL268W_sgtR_URA3"

# Update Scale values (C2, C4) from 25nm to 100nm
$ws.Range("C2").Value = "100nm"
$ws.Range("C4").Value = "100nm"

# Update selection to A6
$ws.Range("A6").Select()
